# Auto-applies the weekly update to the Nispero (Vega Modelo de Temuco) price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44901
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 2500
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 45219
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 35000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 35000
$ws.Range("R3").Value = 'Provincia de Los Andes'
$ws.Range("S3").Value = 3500

# Row 4
$ws.Range("D4").Value = 44519
$ws.Range("M4").Value = 30
$ws.Range("R4").Value = 'Provincia de Quillota'

# Row 5
$ws.Range("D5").Value = 44488
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/bandeja 5 kilos'
$ws.Range("R5").Value = 'La Ligua'
$ws.Range("S5").Value = 2400
$ws.Range("T5").Value = 5

# Row 6
$ws.Range("D6").Value = 44859
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = '$/bandeja 5 kilos'
$ws.Range("S6").Value = 4000
$ws.Range("T6").Value = 5

# Row 7
$ws.Range("D7").Value = 44515
$ws.Range("M7").Value = 80
$ws.Range("R7").Value = 'Provincia de Los Andes'

# Row 8
$ws.Range("D8").Value = 44858
$ws.Range("M8").Value = 90
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("R8").Value = 'Provincia de Quillota'
$ws.Range("S8").Value = 4000

# Row 9
$ws.Range("D9").Value = 44483
$ws.Range("M9").Value = 35
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = '$/bandeja 5 kilos'
$ws.Range("S9").Value = 2000
$ws.Range("T9").Value = 5

# Row 10
$ws.Range("D10").Value = 44166
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'La Ligua'
$ws.Range("S10").Value = 667
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44879
$ws.Range("M11").Value = 25
$ws.Range("N11").Value = 30000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 30000
$ws.Range("R11").Value = 'Provincia de Quillota'
$ws.Range("S11").Value = 3000

# Row 12
$ws.Range("D12").Value = 44496
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 28000
$ws.Range("O12").Value = 28000
$ws.Range("P12").Value = 28000
$ws.Range("Q12").Value = '$/bandeja 10 kilos'
$ws.Range("R12").Value = 'Provincia de Quillota'
$ws.Range("S12").Value = 2800
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("D13").Value = 44466
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 11000
$ws.Range("O13").Value = 11000
$ws.Range("P13").Value = 11000
$ws.Range("R13").Value = 'La Ligua'
$ws.Range("S13").Value = 2200

# Row 15
$ws.Range("D15").Value = 44868
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = '$/bandeja 5 kilos'
$ws.Range("T15").Value = 5

# Row 16
$ws.Range("D16").Value = 45222
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 28000
$ws.Range("O16").Value = 28000
$ws.Range("P16").Value = 28000
$ws.Range("S16").Value = 2800

# Row 17
$ws.Range("D17").Value = 44889
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 30000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 30000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = 'Provincia de Quillota'
$ws.Range("S17").Value = 3000
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("D18").Value = 44511
$ws.Range("M18").Value = 45

# Row 19
$ws.Range("D19").Value = 44511
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 3200
$ws.Range("O19").Value = 3200
$ws.Range("P19").Value = 3200
$ws.Range("S19").Value = 320

# Row 20
$ws.Range("D20").Value = 44503
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 28000
$ws.Range("O20").Value = 28000
$ws.Range("P20").Value = 28000
$ws.Range("S20").Value = 2800

# Row 21
$ws.Range("D21").Value = 44874
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 25000
$ws.Range("O21").Value = 25000
$ws.Range("P21").Value = 25000
$ws.Range("Q21").Value = '$/bandeja 10 kilos'
$ws.Range("S21").Value = 2500
$ws.Range("T21").Value = 10

# Row 22
$ws.Range("D22").Value = 44921
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("Q22").Value = '$/bandeja 7 kilos'
$ws.Range("S22").Value = 2143
$ws.Range("T22").Value = 7

# Row 23
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = 'Vega Modelo de Temuco'
$ws.Range("C23").Value = 'La Araucanía'
$ws.Range("D23").Value = 45224
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 'Fruta'
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = 'Frutos de pepita'
$ws.Range("I23").Value = 100104004
$ws.Range("J23").Value = 'Níspero'
$ws.Range("K23").Value = 'Californiana(o)'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 40
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("R23").Value = 'Provincia de Los Andes'
$ws.Range("S23").Value = 3000
$ws.Range("T23").Value = 10

